$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 175; existing rows 175-188 shift down to 176-189.
$ws.Rows(175).Insert()

# Populate the newly inserted row 175 with the new weekly price observation.
$ws.Cells.Item(175, 1).Value2 = 1
$ws.Cells.Item(175, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(175, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(175, 4).Value2 = 45021
$ws.Cells.Item(175, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(175, 5).Value2 = 15
$ws.Cells.Item(175, 6).Value2 = 100114001
$ws.Cells.Item(175, 7).Value2 = "Papa"
$ws.Cells.Item(175, 8).Value2 = "Asterix"
$ws.Cells.Item(175, 9).Value2 = "1a (cosecha)"
$ws.Cells.Item(175, 10).Value2 = 1000
$ws.Cells.Item(175, 11).Value2 = 13000
$ws.Cells.Item(175, 12).Value2 = 14000
$ws.Cells.Item(175, 13).Value2 = 13500
$ws.Cells.Item(175, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(175, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(175, 16).Value2 = 540
$ws.Cells.Item(175, 17).Value2 = 25
$ws.Cells.Item(175, 18).Value2 = "Hortaliza"
